$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(36).Insert()

$ws.Range("A36").Value = 8
$ws.Range("B36").Value = "Terminal La Palmera de La Serena"
$ws.Range("C36").Value = "Coquimbo"
$ws.Range("D36").Value = 45044
$ws.Range("E36").Value = 4
$ws.Range("F36").Value = 100114007
$ws.Range("G36").Value = "Jengibre"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 400
$ws.Range("K36").Value = 16500
$ws.Range("L36").Value = 17000
$ws.Range("M36").Value = 16750
$ws.Range("N36").Value = "$/caja 13 kilos"
$ws.Range("O36").Value = "Perú"
$ws.Range("P36").Value = 1288
$ws.Range("Q36").Value = 13
$ws.Range("R36").Value = "Hortaliza"
